$wb = $excel.ActiveWorkbook

$sheetName = "ROW50-FE-LIFTER"
$ws = $wb.Worksheets.Item($sheetName)
$rows = @(
    @{ row = 386; A = 45721.73015877315; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"; E = "0x14"; I = 20 },
    @{ row = 387; A = 45721.73018203703; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"; E = "0x14"; I = 20 },
    @{ row = 388; A = 45721.7302053125; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"; E = "0x14"; I = 20 },
    @{ row = 389; A = 45722.23047579861; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"; E = "0x14"; I = 20 },
    @{ row = 390; A = 45722.23049802084; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"; E = "0x14"; I = 20 },
    @{ row = 391; A = 45722.23052140047; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"; E = "0x14"; I = 20 },
    @{ row = 392; A = 45723.19127907408; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"; E = "0x14"; I = 20 },
    @{ row = 393; A = 45723.19130241898; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"; E = "0x14"; I = 20 },
    @{ row = 394; A = 45723.19132570602; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"; E = "0x14"; I = 20 }
)
foreach ($item in $rows) {
    $r = $item.row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 2).Value = "0x01,0x90"
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = "0x01,0x90,"
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = 400
    $ws.Cells.Item($r, 7).Value = [double]"5.686312626471138e+23"
    $ws.Cells.Item($r, 8).Value = 400
    $ws.Cells.Item($r, 9).Value = $item.I
}

$sheetName = "ROW50-MID-LIFTER"
$ws = $wb.Worksheets.Item($sheetName)
$rows = @(
    @{ row = 501; A = 45721.72839921296; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"; E = "0x19"; I = 25 },
    @{ row = 502; A = 45721.7284225; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"; E = "0x19"; I = 25 },
    @{ row = 503; A = 45721.72844564815; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"; E = "0x19"; I = 25 },
    @{ row = 504; A = 45722.22854285879; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"; E = "0x19"; I = 25 },
    @{ row = 505; A = 45722.22856443287; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"; E = "0x19"; I = 25 },
    @{ row = 506; A = 45722.22858758102; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"; E = "0x19"; I = 25 },
    @{ row = 507; A = 45722.72868413194; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"; E = "0x19"; I = 25 },
    @{ row = 508; A = 45722.72870657407; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"; E = "0x19"; I = 25 },
    @{ row = 509; A = 45722.72872986111; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"; E = "0x19"; I = 25 },
    @{ row = 510; A = 45723.22882704861; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"; E = "0x19"; I = 25 },
    @{ row = 511; A = 45723.22884876157; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"; E = "0x19"; I = 25 },
    @{ row = 512; A = 45723.22887202547; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"; E = "0x19"; I = 25 }
)
foreach ($item in $rows) {
    $r = $item.row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 2).Value = "0x01,0x90"
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = "0x01,0x90,"
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = 400
    $ws.Cells.Item($r, 7).Value = [double]"5.686312626471138e+23"
    $ws.Cells.Item($r, 8).Value = 400
    $ws.Cells.Item($r, 9).Value = $item.I
}

$sheetName = "ROW11-FE-LIFTER"
$ws = $wb.Worksheets.Item($sheetName)
$rows = @(
    @{ row = 425; A = 45721.72979140046; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"; E = "0x13"; I = 19 },
    @{ row = 426; A = 45721.72981465278; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"; E = "0x13"; I = 19 },
    @{ row = 427; A = 45721.72983799769; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"; E = "0x13"; I = 19 },
    @{ row = 428; A = 45722.23010865741; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"; E = "0x13"; I = 19 },
    @{ row = 429; A = 45722.23013063658; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"; E = "0x13"; I = 19 },
    @{ row = 430; A = 45722.23015388889; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"; E = "0x13"; I = 19 },
    @{ row = 431; A = 45723.19126707176; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"; E = "0x13"; I = 19 },
    @{ row = 432; A = 45723.19129023148; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"; E = "0x13"; I = 19 },
    @{ row = 433; A = 45723.19131362269; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"; E = "0x13"; I = 19 }
)
foreach ($item in $rows) {
    $r = $item.row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 2).Value = "0x01,0x90"
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = "0x01,0x90,"
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = 400
    $ws.Cells.Item($r, 7).Value = [double]"5.686312626471138e+23"
    $ws.Cells.Item($r, 8).Value = 400
    $ws.Cells.Item($r, 9).Value = $item.I
}

$sheetName = "ROW11-MID-LIFTER"
$ws = $wb.Worksheets.Item($sheetName)
$rows = @(
    @{ row = 362; A = 45721.73010056713; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"; E = "0x9"; I = 9 },
    @{ row = 363; A = 45721.73012371528; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"; E = "0x9"; I = 9 },
    @{ row = 364; A = 45721.73014709491; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"; E = "0x9"; I = 9 },
    @{ row = 365; A = 45722.23024466435; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"; E = "0x9"; I = 9 },
    @{ row = 366; A = 45722.23026591435; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"; E = "0x9"; I = 9 },
    @{ row = 367; A = 45722.23028918981; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"; E = "0x9"; I = 9 },
    @{ row = 368; A = 45723.19114016204; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"; E = "0x9"; I = 9 },
    @{ row = 369; A = 45723.19116321759; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"; E = "0x9"; I = 9 },
    @{ row = 370; A = 45723.19118659722; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"; E = "0x9"; I = 9 }
)
foreach ($item in $rows) {
    $r = $item.row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 2).Value = "0x01,0x90"
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = "0x01,0x90,"
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = 400
    $ws.Cells.Item($r, 7).Value = [double]"5.686312626471138e+23"
    $ws.Cells.Item($r, 8).Value = 400
    $ws.Cells.Item($r, 9).Value = $item.I
}
